$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.601.19'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.940.22'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9982'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.89'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9977'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4847'
$ws.Range('E7').Value = '  +2.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2907'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06811'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '112.65'
$ws.Range('E10').Value = '  +6.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.50'
$ws.Range('E11').Value = '  +5.80%  '
$ws.Range('D12').Value = '1.927.95'
$ws.Range('E12').Value = '  +0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.496'
$ws.Range('E13').Value = '  +2.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07579'
$ws.Range('E14').Value = '  -1.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6836'
$ws.Range('E15').Value = '  +1.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '297.05'
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('D17').Value = '30.556.57'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007683'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.11'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.579'
$ws.Range('E20').Value = '  -0.31%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9986'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '2.170.50'
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9986'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.492'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.507'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.23'
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.55'
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.132'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1075'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.446'
$ws.Range('E30').Value = '  +2.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.166'
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.084'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04988'
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7398'
$ws.Range('E34').Value = '  +0.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.147'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02040'
$ws.Range('E36').Value = '  -1.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.697'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.032'
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '109.91'
$ws.Range('E40').Value = '  -1.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4473'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8697'
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.869'
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.80'
$ws.Range('E44').Value = '  +2.64%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '49.54'
$ws.Range('E46').Value = '  +3.19%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.293'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.330'
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1236'
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.2523'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.02'
$ws.Range('E51').Value = '  -0.84%  '
